$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Klay Thompson -> Jalen Green
$ws.Range("A3").Value = "Jalen Green"
$ws.Range("B3").Value = "PG,SG"
$ws.Range("C3").Value = "Houston Rockets"

# Row 11: Nikola Jokić -> Nikola Jokic (diacritic removed)
$ws.Range("A11").Value = "Nikola Jokic"

# Row 12: Jaylen Brown -> Clint Capela
$ws.Range("A12").Value = "Clint Capela"
$ws.Range("B12").Value = "C"
$ws.Range("C12").Value = "Atlanta Hawks"

# Row 13: Clint Capela -> Jakob Poltl
$ws.Range("A13").Value = "Jakob Poltl"
$ws.Range("B13").Value = "C"
$ws.Range("C13").Value = "Toronto Raptors"

# Row 14: Jalen Green -> Klay Thompson
$ws.Range("A14").Value = "Klay Thompson"
$ws.Range("B14").Value = "SG,SF"
$ws.Range("C14").Value = "Dallas Mavericks"

# Row 15: Jakob Pöltl -> Jaylen Brown
$ws.Range("A15").Value = "Jaylen Brown"
$ws.Range("B15").Value = "SG,SF"
$ws.Range("C15").Value = "Boston Celtics"
